$wb = $excel.ActiveWorkbook

# Worksheets (1-based, in physical/tab order):
#   1 = first
#   2 = GeneralTaxRateWeekly      -> GeneralTaxRateMonthly
#   3 = ProcessPayrollForWeeklyTax -> ProcessPayrollForMonthlyTax
#   4 = TestReports
$wsFirst = $wb.Worksheets.Item(1)
$wsGen   = $wb.Worksheets.Item(2)
$wsProc  = $wb.Worksheets.Item(3)
$wsTest  = $wb.Worksheets.Item(4)

# Rename the "Weekly" sheets to "Monthly" - this is the core of the edit:
# the input sheet was repurposed from a weekly payroll scenario to a
# monthly one.
$wsGen.Name  = "GeneralTaxRateMonthly"
$wsProc.Name = "ProcessPayrollForMonthlyTax"

# Update the matching descriptive cell values on the "first" sheet so the
# test-case table still references the correct (renamed) worksheet names.
$wsFirst.Range("A3").Value = "GeneralTaxRateMonthly"
$wsFirst.Range("A4").Value = "ProcessPayrollForMonthlyTax"

# The "DO NOT TOUCH AUTOMATION EMP 107" marker (shared across three
# worksheets) becomes "... EMP 105".
$wsGen.Range("A2").Value  = "DO NOT TOUCH AUTOMATION EMP 105"
$wsProc.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsTest.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

# Restore each sheet's last-used cell selection, and leave
# "GeneralTaxRateMonthly" as the active tab (select it last).
$wsFirst.Range("F5").Select()
$wsProc.Range("H12").Select()
$wsTest.Range("B6").Select()
$wsGen.Range("E6").Select()
